$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(118, 8).Value = 7916
$ws.Cells.Item(118, 9).Value = 9570
$ws.Cells.Item(118, 10).Value = 1300
$ws.Cells.Item(118, 11).Value = 28710
$ws.Cells.Item(118, 12).Value = 3900
$ws.Cells.Item(118, 13).Value = -27053
$ws.Cells.Item(118, 14).Value = -7214

$ws.Cells.Item(132, 8).Value = 7149218
$ws.Cells.Item(132, 9).Value = 8071581.5
$ws.Cells.Item(132, 10).Value = 899.75
$ws.Cells.Item(132, 11).Value = 24214744.5
$ws.Cells.Item(132, 12).Value = 2699.25
$ws.Cells.Item(132, 13).Value = -24212214.5
$ws.Cells.Item(132, 14).Value = -7759.25

$ws.Cells.Item(134, 8).Value = 0
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 14).ClearContents()  # N134 cell removed entirely per diff

$ws.Cells.Item(135, 8).Value = 914.6842
$ws.Cells.Item(135, 9).Value = 592.6445
$ws.Cells.Item(135, 10).Value = 2122.3333
$ws.Cells.Item(135, 11).Value = 5333.8005
$ws.Cells.Item(135, 12).Value = 19100.9997
$ws.Cells.Item(135, 13).Value = -2798.8005
$ws.Cells.Item(135, 14).Value = -24170.9997

$ws.Cells.Item(138, 8).Value = 1789.6
$ws.Cells.Item(138, 9).Value = 1153.6666
$ws.Cells.Item(138, 10).Value = 3424.8572
$ws.Cells.Item(138, 11).Value = 3460.9998
$ws.Cells.Item(138, 12).Value = 10274.5716
$ws.Cells.Item(138, 13).Value = 1679.0002
$ws.Cells.Item(138, 14).Value = -20554.5716

$ws.Cells.Item(141, 8).Value = 1297.3784
$ws.Cells.Item(141, 9).Value = 1128.1029
$ws.Cells.Item(141, 11).Value = 3384.3087
$ws.Cells.Item(141, 13).Value = 1795.6913

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 20641.895
$ws.Cells.Item(32, 9).Value = 4573.5767
$ws.Cells.Item(32, 11).Value = 4573.5767
$ws.Cells.Item(32, 13).Value = -4286.5767

$ws.Cells.Item(61, 8).Value = 1367.697
$ws.Cells.Item(61, 9).Value = 1086.4814
$ws.Cells.Item(61, 10).Value = 2633.1667
$ws.Cells.Item(61, 11).Value = 1086.4814
$ws.Cells.Item(61, 12).Value = 2633.1667
$ws.Cells.Item(61, 13).Value = -874.4813999999999
$ws.Cells.Item(61, 14).Value = -3057.1667

$ws.Cells.Item(74, 8).Value = 812.2222
$ws.Cells.Item(74, 9).Value = 656.5
$ws.Cells.Item(74, 11).Value = 656.5
$ws.Cells.Item(74, 13).Value = 217.5

$ws.Cells.Item(77, 8).Value = 812.2222
$ws.Cells.Item(77, 9).Value = 656.5
$ws.Cells.Item(77, 11).Value = 3282.5
$ws.Cells.Item(77, 13).Value = 1085.5

$ws.Cells.Item(132, 8).Value = 8354.134
$ws.Cells.Item(132, 9).Value = 10127.272
$ws.Cells.Item(132, 10).Value = 3478
$ws.Cells.Item(132, 11).Value = 30381.816
$ws.Cells.Item(132, 12).Value = 10434
$ws.Cells.Item(132, 13).Value = -27851.816
$ws.Cells.Item(132, 14).Value = -15494

$ws.Cells.Item(136, 8).Value = 1367.697
$ws.Cells.Item(136, 9).Value = 1086.4814
$ws.Cells.Item(136, 10).Value = 2633.1667
$ws.Cells.Item(136, 11).Value = 3259.4442
$ws.Cells.Item(136, 12).Value = 7899.500100000001
$ws.Cells.Item(136, 13).Value = -709.4441999999999
$ws.Cells.Item(136, 14).Value = -12999.5001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2461.4634
$ws.Cells.Item(134, 9).Value = 2219.7222
$ws.Cells.Item(134, 10).Value = 4202
$ws.Cells.Item(134, 11).Value = 6659.1666
$ws.Cells.Item(134, 12).Value = 12606
$ws.Cells.Item(134, 13).Value = -4124.1666
$ws.Cells.Item(134, 14).Value = -17676

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 43235.855
$ws.Cells.Item(31, 9).Value = 2171.5
$ws.Cells.Item(31, 10).Value = 77816.37
$ws.Cells.Item(31, 11).Value = 2171.5
$ws.Cells.Item(31, 12).Value = 77816.37
$ws.Cells.Item(31, 13).Value = -1876.5
$ws.Cells.Item(31, 14).Value = -78406.37

$ws.Cells.Item(34, 8).Value = 43235.855
$ws.Cells.Item(34, 9).Value = 2171.5
$ws.Cells.Item(34, 10).Value = 77816.37
$ws.Cells.Item(34, 11).Value = 2171.5
$ws.Cells.Item(34, 12).Value = 77816.37
$ws.Cells.Item(34, 13).Value = -1969.5
$ws.Cells.Item(34, 14).Value = -78220.37

$ws.Cells.Item(58, 8).Value = 1196.4259
$ws.Cells.Item(58, 9).Value = 1059.591
$ws.Cells.Item(58, 10).Value = 1798.5
$ws.Cells.Item(58, 11).Value = 1059.591
$ws.Cells.Item(58, 12).Value = 1798.5
$ws.Cells.Item(58, 13).Value = -856.5909999999999
$ws.Cells.Item(58, 14).Value = -2204.5

$ws.Cells.Item(132, 8).Value = 3723.0557
$ws.Cells.Item(132, 9).Value = 3673.6924
$ws.Cells.Item(132, 10).Value = 3851.4
$ws.Cells.Item(132, 11).Value = 11021.0772
$ws.Cells.Item(132, 12).Value = 11554.2
$ws.Cells.Item(132, 13).Value = -8491.0772
$ws.Cells.Item(132, 14).Value = -16614.2

$ws.Cells.Item(134, 8).Value = 1135.0714
$ws.Cells.Item(134, 9).Value = 1189.0555
$ws.Cells.Item(134, 10).Value = 1037.9
$ws.Cells.Item(134, 11).Value = 3567.1665
$ws.Cells.Item(134, 12).Value = 3113.7
$ws.Cells.Item(134, 13).Value = -1032.1665
$ws.Cells.Item(134, 14).Value = -8183.700000000001

$ws.Cells.Item(136, 8).Value = 1196.4259
$ws.Cells.Item(136, 9).Value = 1059.591
$ws.Cells.Item(136, 10).Value = 1798.5
$ws.Cells.Item(136, 11).Value = 3178.773
$ws.Cells.Item(136, 12).Value = 5395.5
$ws.Cells.Item(136, 13).Value = -628.7729999999997
$ws.Cells.Item(136, 14).Value = -10495.5

$ws.Cells.Item(138, 8).Value = 55361.25
$ws.Cells.Item(138, 10).Value = 55361.25
$ws.Cells.Item(138, 12).Value = 55361.25
$ws.Cells.Item(138, 14).Value = -65641.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 962.8372000000001
$ws.Cells.Item(5, 9).Value = 526.2174
$ws.Cells.Item(5, 10).Value = 1464.95
$ws.Cells.Item(5, 11).Value = 1578.6522
$ws.Cells.Item(5, 12).Value = 4394.85
$ws.Cells.Item(5, 13).Value = -1466.6522
$ws.Cells.Item(5, 14).Value = -4618.85

$ws.Cells.Item(122, 8).Value = 505
$ws.Cells.Item(122, 10).Value = 440
$ws.Cells.Item(122, 12).Value = 3960
$ws.Cells.Item(122, 14).Value = -8860

$ws.Cells.Item(135, 8).Value = 962.8372000000001
$ws.Cells.Item(135, 9).Value = 526.2174
$ws.Cells.Item(135, 10).Value = 1464.95
$ws.Cells.Item(135, 11).Value = 4735.9566
$ws.Cells.Item(135, 12).Value = 13184.55
$ws.Cells.Item(135, 13).Value = -2200.9566
$ws.Cells.Item(135, 14).Value = -18254.55

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3261.0244
$ws.Cells.Item(132, 9).Value = 3133.7742
$ws.Cells.Item(132, 10).Value = 3655.5
$ws.Cells.Item(132, 11).Value = 9401.3226
$ws.Cells.Item(132, 12).Value = 10966.5
$ws.Cells.Item(132, 13).Value = -6871.3226
$ws.Cells.Item(132, 14).Value = -16026.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1125177.5
$ws.Cells.Item(46, 9).Value = 429.33334
$ws.Cells.Item(46, 10).Value = 1687551.5
$ws.Cells.Item(46, 11).Value = 429.33334
$ws.Cells.Item(46, 12).Value = 1687551.5
$ws.Cells.Item(46, 13).Value = -241.33334
$ws.Cells.Item(46, 14).Value = -1687927.5

$ws.Cells.Item(132, 8).Value = 3463.074
$ws.Cells.Item(132, 9).Value = 4372.1113
$ws.Cells.Item(132, 10).Value = 1645
$ws.Cells.Item(132, 11).Value = 13116.3339
$ws.Cells.Item(132, 12).Value = 4935
$ws.Cells.Item(132, 13).Value = -10586.3339
$ws.Cells.Item(132, 14).Value = -9995

$ws.Cells.Item(136, 8).Value = 1316.921
$ws.Cells.Item(136, 9).Value = 1104.0857
$ws.Cells.Item(136, 10).Value = 3800
$ws.Cells.Item(136, 11).Value = 3312.2571
$ws.Cells.Item(136, 12).Value = 11400
$ws.Cells.Item(136, 13).Value = -762.2571000000003
$ws.Cells.Item(136, 14).Value = -16500

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 143942.86
$ws.Cells.Item(100, 9).Value = 251175
$ws.Cells.Item(100, 10).Value = 966.6667
$ws.Cells.Item(100, 11).Value = 502350
$ws.Cells.Item(100, 12).Value = 1933.3334
$ws.Cells.Item(100, 13).Value = -501809
$ws.Cells.Item(100, 14).Value = -3015.3334

$ws.Cells.Item(136, 8).Value = 555.1525
$ws.Cells.Item(136, 9).Value = 380.18182
$ws.Cells.Item(136, 10).Value = 1068.4
$ws.Cells.Item(136, 11).Value = 1140.54546
$ws.Cells.Item(136, 12).Value = 3205.2
$ws.Cells.Item(136, 13).Value = 1409.45454
$ws.Cells.Item(136, 14).Value = -8305.200000000001
